# Apply "most recent excel files" update to Test.xlsx.
# New equipment/vessel/voyage records are added, which pushes the
# previously-existing row 4 record (RFCU4114542 / EVER LEGACY) down,
# and appends several additional rows through row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CCLU4715365", "EVER SMILE",  "1081E", "9077904371-01", "", "COSU6204572830"),
    @("FSCU5038306", "EVER SMILE",  "1081E", "9077904373-01", "", "COSU6204572830"),
    @("RFCU4114542", "EVER LEGACY", "034E",  "9072904799-01", "", "COSU6203869480"),
    @("CBHU6409087", "EVER SMILE",  "1081E", "9077904368-01", "", "COSU6204572830"),
    @("HESU4031448", "EVER LOVELY", "0850E", "9071905305-01", "", "COSU6211257620"),
    @("CCLU4829579", "EVER SMILE",  "1081E", "9077904372-01", "", "COSU6204572830"),
    @("HESU4027089", "NAVARINO",    "0848E", "9075904134-01", "", "COSU6203956310"),
    @("HESU4031406", "EVER LOVELY", "0850E", "9071905304-01", "", "COSU6211257620")
)

$rowIndex = 4
foreach ($row in $data) {
    $ws.Range("A$rowIndex").Value = $row[0]
    $ws.Range("B$rowIndex").Value = $row[1]
    $ws.Range("C$rowIndex").Value = $row[2]
    $ws.Range("D$rowIndex").Value = $row[3]
    $ws.Range("F$rowIndex").Value = $row[5]
    $rowIndex++
}

Write-Host "Applied updates through row $($rowIndex - 1)"
